# "Colocando header nos gráficos"
# Adds a header label in column A (row 1) for each data table used by the
# charts, fixes missing Portuguese accents in several row labels, removes
# the "Teto" row from the emissions sheet, and updates the cost figures /
# header on the "Custo Total" sheet.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell {
    param($ws, [string]$formatSourceAddress, [string]$targetAddress, [string]$text)

    # Copy number/font/border/alignment formatting from a neighboring header
    # cell (e.g. B1) onto the new A1 header cell, then set its text.
    $ws.Range($formatSourceAddress).Copy() | Out-Null
    $ws.Range($targetAddress).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($targetAddress).Value = $text
}

function Clear-LabelFormat {
    param($ws, [string]$address, [string]$text)

    $ws.Range($address).ClearFormats() | Out-Null
    if ($null -ne $text) {
        $ws.Range($address).Value = $text
    }
}

function Set-TextValue {
    param($ws, [string]$address, [string]$text, [string]$styleSourceAddress)

    # Force a number-looking string (e.g. "2015") to be stored as text rather
    # than being auto-converted to a numeric cell, then restore the header
    # formatting (bold/border/center) that the plain assignment would disturb.
    $ws.Range($address).NumberFormat = "@"
    $ws.Range($address).Value = $text
    if ($styleSourceAddress) {
        $ws.Range($styleSourceAddress).Copy() | Out-Null
        $ws.Range($address).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    }
}

# ---------------------------------------------------------------------
# Sheets 1-4 share the same row layout (Fonte/Tecnologia breakdown table)
# ---------------------------------------------------------------------
$sourceSheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($sheetName in $sourceSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    Set-HeaderCell $ws "B1" "A1" "Fonte/Tecnologia"

    Clear-LabelFormat $ws "A2" $null                  # Hidro (unchanged text)
    Clear-LabelFormat $ws "A3" "Gás Natural"           # Gas Natural -> Gás Natural
    Clear-LabelFormat $ws "A4" "Carvão"                # Carvao -> Carvão
    Clear-LabelFormat $ws "A5" $null                   # Nuclear (unchanged text)
    Clear-LabelFormat $ws "A6" "Óleos Comb"            # Oleos Comb -> Óleos Comb
    Clear-LabelFormat $ws "A7" $null                   # Biomassa (unchanged text)
    Clear-LabelFormat $ws "A8" "Eólica"                # Eolica -> Eólica
    Clear-LabelFormat $ws "A9" $null                   # Solar (unchanged text)
    Clear-LabelFormat $ws "A10" $null                  # Outros (unchanged text)
    Clear-LabelFormat $ws "A11" "Pot. Compl."          # Pot Compl -> Pot. Compl.
    Clear-LabelFormat $ws "A12" $null                  # GD (unchanged text)
}

# ---------------------------------------------------------------------
# Sheet 5 - Emissoes Totais (MtCO2eq): add header, fix labels, drop "Teto"
# ---------------------------------------------------------------------
$wsEmissoes = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

Set-HeaderCell $wsEmissoes "B1" "A1" "Período"

Clear-LabelFormat $wsEmissoes "A2" "P.Médio"           # P Medio -> P.Médio
Clear-LabelFormat $wsEmissoes "A3" "P.Crítico"         # P Critico -> P.Crítico

# Remove the "Teto" row entirely (row 4)
$wsEmissoes.Rows("4:4").Delete() | Out-Null

# ---------------------------------------------------------------------
# Sheet 6 - Custo Total (bilhões de R$): add header, fix labels & values
# ---------------------------------------------------------------------
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

Set-HeaderCell $wsCusto "B1" "A1" "Tipo Expansão"
Set-TextValue $wsCusto "B1" "2015" "A2"

Clear-LabelFormat $wsCusto "A2" "Expansão Centralizada"  # Expansao Centralizada -> Expansão Centralizada
$wsCusto.Range("B2").Value = 573

Clear-LabelFormat $wsCusto "A3" "Expansão por GD"        # Expansao por GD -> Expansão por GD
$wsCusto.Range("B3").Value = 99
